$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.366.89"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.077.71"
$ws.Range("E3").Value = "  +5.17%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "236.37"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "57.11"
$ws.Range("E8").Value = "  +5.69%  "
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").Value = "58.21"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "0.0761"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "2.383.05"
$ws.Range("E13").Value = "  +5.25%  "
$ws.Range("D14").Value = "14.57"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "21.18"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "0.776"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "2.081.87"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("D19").Value = "37.544.00"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "5.93"
$ws.Range("E20").Value = "  +20.00%  "
$ws.Range("D21").Value = "68.44"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "224.19"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("D27").Value = "163.72"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("D29").Value = "0.132"
$ws.Range("E29").Value = "  +7.54%  "
$ws.Range("D30").Value = "19.27"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  +6.90%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "4.47"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "0.0620"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +9.63%  "
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "5.97"
$ws.Range("E38").Value = "  +15.40%  "
$ws.Range("D39").Value = "3.33"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  +23.49%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0957"
$ws.Range("E43").Value = "  +8.07%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.474.56"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("D45").Value = "94.98"
$ws.Range("E45").Value = "  +8.96%  "
$ws.Range("D46").Value = "0.0209"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "16.09"
$ws.Range("E48").Value = "  +5.84%  "
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").Value = "7.16"
$ws.Range("E50").Value = "  +8.56%  "
$ws.Range("D51").Value = "2.95"
$ws.Range("E51").Value = "  +3.42%  "
